# Applies the diff: appends "." after the YouTube hyperlink, then adds a
# block of new paragraphs (a page-break note, spacing, and a styled note
# about the VARCHAR "ALIAS" column) right before the document's final
# (empty) paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: insert a paragraph-level OOXML fragment immediately before the
# current last paragraph of the document (i.e. right before the body's
# trailing empty paragraph / sectPr). Using the pkg:package-wrapped
# WordOpenXML form lets Word parse real <w:p> structure (pPr/rPr, rStyle,
# lastRenderedPageBreak, ...) instead of just inserting plain text runs,
# and it does so without leaving a stray empty <w:r/> behind the way
# InsertParagraphAfter()/Range.Style do in this host.
# ---------------------------------------------------------------------
function Insert-BodyParagraphXml {
    param([string]$InnerBodyXml)

    $target = $d.Paragraphs.Last.Range
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $InnerBodyXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

# 1) Find the paragraph that holds the YouTube hyperlink and append a
#    plain "." run right after it (still inside the same paragraph).
$findRng = $d.Content
$findRng.Find.Execute("https://www.youtube.com/watch?v=86zyFtd4FkE&t=7s", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$precedingRng = $d.Range(0, $findRng.End)
$hyperlinkParaIndex = $precedingRng.Paragraphs.Count
$hyperlinkPara = $d.Paragraphs.Item($hyperlinkParaIndex)
$hyperlinkPara.Range.InsertAfter(".")

# 2) New paragraph: the "IMPORTANTE..." note, which starts on a new page
#    in the original render (carries a lastRenderedPageBreak marker).
Insert-BodyParagraphXml('<w:p><w:r><w:lastRenderedPageBreak/><w:t>IMPORTANTE&#161;! AL MOMENTO DE REGISTRAR UNA VENTA AVERIGUAR COMO LLEVAR UN HISTORIAL DE LOS PRODUCTOS VENDIDOS EN CADA VENTA.</w:t></w:r></w:p>')

# 3) A blank spacer paragraph.
Insert-BodyParagraphXml('<w:p/>')

# 4) New paragraph: the "ADEMAS..." note, whose paragraph mark carries the
#    Hyperlink character style.
Insert-BodyParagraphXml('<w:p><w:pPr><w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr></w:pPr><w:r><w:t>ADEMAS DE CAMBIAR LA CANTIDAD DE EL VARCHAR ID DE LA TABLA VENTAS &#8220;ALIAS&#8221; a 20.</w:t></w:r></w:p>')

# 5) Two empty paragraphs whose marks still carry the Hyperlink character
#    style (matches the pasted-blank-lines pattern after the note above).
Insert-BodyParagraphXml('<w:p><w:pPr><w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr></w:pPr></w:p>')
Insert-BodyParagraphXml('<w:p><w:pPr><w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr></w:pPr></w:p>')

# 6) A final empty paragraph whose mark keeps the hyperlink's direct
#    color/underline formatting (style link dropped, direct formatting
#    remains) before the document's own trailing empty paragraph.
Insert-BodyParagraphXml('<w:p><w:pPr><w:rPr><w:color w:val="0563C1" w:themeColor="hyperlink"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>')

Write-Output "done; paragraph count = $($d.Paragraphs.Count)"
